$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of column letter -> new value, applied to rows 2 and 3
$updates = @{
    "D"  = 0.00735
    "E"  = -0.0118
    "G"  = 0.09896355733868271
    "H"  = 0.09896355733868271
    "I"  = 0.08993647609495151
    "J"  = 0.07424330910167347
    "K"  = 44.2
    "L"  = 0.1477766633233032
    "M"  = 30.3
    "N"  = 0.06614276358873608
    "O"  = 0.6855203619909502
    "P"  = 30.3
    "Q"  = 0.06614276358873608
    "R"  = 0.6855203619909502
    "U"  = 9.640000000000001
    "V"  = 0.02104344029687841
    "W"  = 0.1705905055962949
    "X"  = 0.04602442524895171
    "Y"  = 0.1245660803473431
    "Z"  = 1.199278267842823
    "AA" = 0.08903838713837424
    "AB" = 0.04596597123634667
    "AC" = 0.04307241590202757
    "AD" = 1.6
    "AE" = 0
    "AF" = 1.6
    "AG" = -8.040000000000001
    "AH" = 0.003480530780944094
    "AI" = 0.006336633663366337
    "AJ" = -0.01786428476203173
    "AK" = -0.0331054928765544
    "AL" = 0.131
    "AM" = 0.131
    "AN" = 0.05776173285198556
    "AO" = 205.3435114503817
    "AP" = -0.2902527075812275
    "AQ" = 205.3435114503817
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("${col}2").Value = $value
    $ws.Range("${col}3").Value = $value
}
